$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row that held only the "5817181 - Valdeir Arantes" value (old row 13),
# which shifts all subsequent rows up by one and naturally carries row heights/labels along.
$ws.Rows.Item(13).Delete()

# After the shift, several value cells need to be updated to their new (shorter) content.
$ws.Range("B10").Value = '5817181 - Valdeir Arantes'
$ws.Range("C10").Value = '5817181 - Valdeir Arantes'

$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'

$ws.Range("B15").Value = '01/01/2019'
$ws.Range("C15").Value = '01/01/2019'

$ws.Range("B18").Value = '5817181 - Valdeir Arantes'
$ws.Range("C18").Value = '5817181 - Valdeir Arantes'

$ws.Range("B19").Value = 'A avaliação será composta por provas, exercícios, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range("C19").Value = 'A avaliação será composta por provas, exercícios, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'

$ws.Range("B20").Value = 'MF≥ 5,0 para aprovação 5,0'
$ws.Range("C20").Value = 'MF≥ 5,0 para aprovação 5,0'

$ws.Range("B21").Value = '(MF+RC)/2 ≥ 5,0 para aprovação, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Range("C21").Value = '(MF+RC)/2 ≥ 5,0 para aprovação, onde RC é uma prova de recuperação a ser aplicada.'
